$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting existing rows 94-118 down to 95-119.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly data entry.
$ws.Range("A94").Value = 7
$ws.Range("B94").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C94").Value = "Ñuble"
$ws.Range("D94").Value = 44463
$ws.Range("E94").Value = 16
$ws.Range("F94").Value = 100112017
$ws.Range("G94").Value = "Apio"
$ws.Range("H94").Value = "Americana (o)"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 160
$ws.Range("K94").Value = 8500
$ws.Range("L94").Value = 9000
$ws.Range("M94").Value = 8750
$ws.Range("N94").Value = "`$/docena de matas"
$ws.Range("O94").Value = "Provincia del Elquí"
$ws.Range("P94").Value = 1458
$ws.Range("Q94").Value = 6
$ws.Range("R94").Value = "Hortaliza"
